$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.264.74'
$ws.Range('E2').Value = '  -2.84%  '
$ws.Range('D3').Value = '1.551.09'
$ws.Range('E3').Value = '  -4.83%  '
$cell = $ws.Range('D5')
$cell.Value = "'206.77"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -3.57%  '
$ws.Range('E6').Value = '  -0.06%  '
$cell = $ws.Range('D7')
$cell.Value = "'0.478"
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -5.22%  '
$cell = $ws.Range('D8')
$cell.Value = "'0.0610"
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('E9').Value = '  -3.23%  '
$cell = $ws.Range('D10')
$cell.Value = "'17.78"
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -4.04%  '
$cell = $ws.Range('D11')
$cell.Value = "'0.0780"
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').Value = '1.763.30'
$ws.Range('E12').Value = '  -4.90%  '
$ws.Range('D13').Value = '1.546.33'
$ws.Range('E13').Value = '  -4.64%  '
$cell = $ws.Range('D14')
$cell.Value = "'3.99"
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -4.66%  '
$cell = $ws.Range('D15')
$cell.Value = "'0.505"
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('D16').Value = '25.217.30'
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('D17').Value = '0.0₃0708'
$ws.Range('E17').Value = '  -4.40%  '
$cell = $ws.Range('D18')
$cell.Value = "'58.67"
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -4.60%  '
$cell = $ws.Range('D19')
$cell.Value = "'1.00"
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.05%  '
$cell = $ws.Range('D20')
$cell.Value = "'186.02"
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -4.02%  '
$cell = $ws.Range('D21')
$cell.Value = "'4.10"
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -3.69%  '
$cell = $ws.Range('D22')
$cell.Value = "'9.25"
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -3.41%  '
$cell = $ws.Range('D23')
$cell.Value = "'5.83"
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -4.25%  '
$ws.Range('B24').Value = 'Stellar'
$ws.Range('C24').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D24')
$cell.Value = "'0.129"
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range('D25')
$cell.Value = "'1.00"
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.11%  '
$cell = $ws.Range('D26')
$cell.Value = "'139.34"
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -3.37%  '
$cell = $ws.Range('D27')
$cell.Value = "'1.63"
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -5.14%  '
$cell = $ws.Range('D28')
$cell.Value = "'14.84"
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -2.97%  '
$cell = $ws.Range('D29')
$cell.Value = "'6.39"
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -5.47%  '
$ws.Range('E30').Value = '  -6.66%  '
$cell = $ws.Range('D31')
$cell.Value = "'0.0463"
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -4.32%  '
$cell = $ws.Range('D32')
$cell.Value = "'3.03"
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -3.62%  '
$cell = $ws.Range('D33')
$cell.Value = "'2.97"
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -4.91%  '
$cell = $ws.Range('D34')
$cell.Value = "'1.46"
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('D36').Value = '1.085.05'
$ws.Range('E36').Value = '  -3.58%  '
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  -2.67%  '
$cell = $ws.Range('D39')
$cell.Value = "'0.493"
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -5.53%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D40')
$cell.Value = "'0.762"
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -10.64%  '
$ws.Range('E41').Value = '  -7.66%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D42')
$cell.Value = "'0.802"
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +4.37%  '
$cell = $ws.Range('D43')
$cell.Value = "'92.78"
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -5.56%  '
$cell = $ws.Range('D44')
$cell.Value = "'5.05"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('D45').Value = '1.678.67'
$ws.Range('E45').Value = '  -4.86%  '
$ws.Range('D46').Value = '0.0₆0108'
$ws.Range('E46').Value = '  -5.31%  '
$cell = $ws.Range('D47')
$cell.Value = "'1.46"
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -1.60%  '
$cell = $ws.Range('D48')
$cell.Value = "'52.33"
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -4.09%  '
$ws.Range('E49').Value = '  -5.69%  '
$ws.Range('E50').Value = '  -0.17%  '
$cell = $ws.Range('D51')
$cell.Value = "'0.404"
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -2.15%  '
